$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price cells hold purely numeric-looking text (e.g. "226.57"); Excel
# auto-converts such strings to numbers on assignment unless the cell is
# pre-formatted as Text, so mark those cells as Text first to preserve the
# original string semantics (matches the source inline-string cells).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "34.161.76"
$ws.Range("D3").Value = "1.788.15"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "226.57"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "0.548"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "31.83"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("D11").Value = "0.0946"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "2.046.05"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("D14").Value = "1.787.54"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "0.624"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").Value = "34.106.34"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").Value = "68.27"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").Value = "247.19"
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").Value = "10.94"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "4.10"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "161.12"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").Value = "7.19"
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("D27").Value = "16.34"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("E31").Value = "  +2.22%  "
$ws.Range("D32").Value = "3.69"
$ws.Range("E32").Value = "  +2.79%  "
$ws.Range("D33").Value = "3.63"
$ws.Range("E33").Value = "  +3.85%  "
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").Value = "1.446.71"
$ws.Range("E35").Value = "  +4.81%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.655"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "2.45"
$ws.Range("E37").Value = "  +9.87%  "
$ws.Range("E38").Value = "  +3.81%  "
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").Value = "80.62"
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("D41").Value = "2.37"
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").Value = "0.924"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").Value = "13.51"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("E45").Value = "  +4.40%  "
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("D49").Value = "1.947.67"
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("D50").Value = "105.96"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("E51").Value = "  -0.01%  "
